# Add new columns (AHV-Nr, Grad, Funktion) to the "Tabelle1" worksheet.
#
# Target layout after the edit:
#   A1=AHV-Nr      B1=Vorname  C1=Nachname  D1=Grad  E1=Funktion      F1=Karte
#   A2=756.9217... B2=Hans     C2=Wurst     D2=Sdt   E2=Stabsassistent F2=1234

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert the two new columns ("Grad"/"Funktion") right after the existing
# "Nachname" column (current column B), pushing "Karte" two slots to the right.
$ws.Columns("C:D").Insert()

# Insert a new column before the existing "Vorname" column (current column A)
# to hold the new "AHV-Nr" column.
$ws.Columns("A:A").Insert()

# Header row
$ws.Range("A1").Value = "AHV-Nr"
$ws.Range("D1").Value = "Grad"
$ws.Range("E1").Value = "Funktion"

# Data row
$ws.Range("A2").Value = "756.9217.0769.85"
$ws.Range("D2").Value = "Sdt"
$ws.Range("E2").Value = "Stabsassistent"

# Size the new columns to fit their content.
$ws.Columns("A:A").AutoFit()
$ws.Columns("E:E").AutoFit()

# Printer / page setup.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore selection to match the saved workbook state.
$ws.Range("G8").Select()
